$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 8 ("grandes regiões" label-only row) first (higher row index first
# so row numbers of rows above are unaffected), then row 5 ("situação do
# domicílio" label-only row). Excel shifts the rows below each deleted row up.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

# Row 2 header cells: B2 and F2 become "total" (same label used in C2).
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"

Write-Output "done"
